# Auto-generated Excel COM-interop script to apply market-price data refresh
# across the Tonberry_Profits leve-profit tables (one table per job sheet).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(41, 8).Value = 284.93332  # H41
$ws.Cells.Item(41, 9).Value = 217.11111  # I41
$ws.Cells.Item(41, 11).Value = 217.11111  # K41
$ws.Cells.Item(41, 13).Value = 222.88889  # M41
$ws.Cells.Item(51, 8).Value = 6268.5713  # H51
$ws.Cells.Item(51, 9).Value = 6576  # I51
$ws.Cells.Item(51, 10).Value = 5500  # J51
$ws.Cells.Item(51, 11).Value = 6576  # K51
$ws.Cells.Item(51, 12).Value = 5500  # L51
$ws.Cells.Item(51, 13).Value = -6092  # M51
$ws.Cells.Item(51, 14).Value = -6468  # N51
$ws.Cells.Item(53, 8).Value = 9254.666999999999  # H53
$ws.Cells.Item(53, 9).Value = 12238.667  # I53
$ws.Cells.Item(53, 10).Value = 302.66666  # J53
$ws.Cells.Item(53, 11).Value = 12238.667  # K53
$ws.Cells.Item(53, 12).Value = 302.66666  # L53
$ws.Cells.Item(53, 13).Value = -11601.667  # M53
$ws.Cells.Item(53, 14).Value = -1576.66666  # N53
$ws.Cells.Item(64, 8).Value = 3666.6667  # H64
$ws.Cells.Item(64, 9).Value = 3000  # I64
$ws.Cells.Item(64, 10).Value = 4000  # J64
$ws.Cells.Item(64, 11).Value = 3000  # K64
$ws.Cells.Item(64, 12).Value = 4000  # L64
$ws.Cells.Item(64, 13).Value = -2752  # M64
$ws.Cells.Item(64, 14).Value = -4496  # N64
$ws.Cells.Item(67, 8).Value = 3666.6667  # H67
$ws.Cells.Item(67, 9).Value = 3000  # I67
$ws.Cells.Item(67, 10).Value = 4000  # J67
$ws.Cells.Item(67, 11).Value = 3000  # K67
$ws.Cells.Item(67, 12).Value = 4000  # L67
$ws.Cells.Item(67, 13).Value = -2142  # M67
$ws.Cells.Item(67, 14).Value = -5716  # N67
$ws.Cells.Item(112, 8).Value = 2653.32  # H112
$ws.Cells.Item(112, 10).Value = 2849.4783  # J112
$ws.Cells.Item(112, 12).Value = 8548.4349  # L112
$ws.Cells.Item(112, 14).Value = -10764.4349  # N112
$ws.Cells.Item(132, 8).Value = 1231.3438  # H132
$ws.Cells.Item(132, 9).Value = 1174.2903  # I132
$ws.Cells.Item(132, 11).Value = 3522.8709  # K132
$ws.Cells.Item(132, 13).Value = -992.8708999999999  # M132
$ws.Cells.Item(141, 8).Value = 4998.75  # H141
$ws.Cells.Item(141, 9).Value = 2995  # I141
$ws.Cells.Item(141, 11).Value = 8985  # K141
$ws.Cells.Item(141, 13).Value = -3805  # M141

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2683.8374  # H32
$ws.Cells.Item(32, 9).Value = 1614.2388  # I32
$ws.Cells.Item(32, 10).Value = 8196.385  # J32
$ws.Cells.Item(32, 11).Value = 1614.2388  # K32
$ws.Cells.Item(32, 12).Value = 8196.385  # L32
$ws.Cells.Item(32, 13).Value = -1327.2388  # M32
$ws.Cells.Item(32, 14).Value = -8770.385  # N32
$ws.Cells.Item(61, 8).Value = 2424.04  # H61
$ws.Cells.Item(61, 9).Value = 1511.6364  # I61
$ws.Cells.Item(61, 11).Value = 1511.6364  # K61
$ws.Cells.Item(61, 13).Value = -1299.6364  # M61
$ws.Cells.Item(97, 8).Value = 550  # H97
$ws.Cells.Item(97, 9).Value = 550  # I97
$ws.Cells.Item(97, 11).Value = 550  # K97
$ws.Cells.Item(97, 13).Value = -54  # M97
$ws.Cells.Item(102, 8).Value = 496.75  # H102
$ws.Cells.Item(102, 9).Value = 496.75  # I102
$ws.Cells.Item(102, 11).Value = 496.75  # K102
$ws.Cells.Item(102, 13).Value = 1125.25  # M102
$ws.Cells.Item(122, 8).Value = 1691.1305  # H122
$ws.Cells.Item(122, 9).Value = 1710.0454  # I122
$ws.Cells.Item(122, 11).Value = 5130.1362  # K122
$ws.Cells.Item(122, 13).Value = -2680.1362  # M122
$ws.Cells.Item(132, 8).Value = 1581.5714  # H132
$ws.Cells.Item(132, 9).Value = 1324.4166  # I132
$ws.Cells.Item(132, 11).Value = 3973.2498  # K132
$ws.Cells.Item(132, 13).Value = -1443.2498  # M132
$ws.Cells.Item(136, 8).Value = 2424.04  # H136
$ws.Cells.Item(136, 9).Value = 1511.6364  # I136
$ws.Cells.Item(136, 11).Value = 4534.9092  # K136
$ws.Cells.Item(136, 13).Value = -1984.9092  # M136

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 2487.3125  # H105
$ws.Cells.Item(105, 9).Value = 2487.3125  # I105
$ws.Cells.Item(105, 11).Value = 2487.3125  # K105
$ws.Cells.Item(105, 13).Value = -740.3125  # M105
$ws.Cells.Item(134, 8).Value = 6482.35  # H134
$ws.Cells.Item(134, 9).Value = 6997.0967  # I134
$ws.Cells.Item(134, 10).Value = 4709.3335  # J134
$ws.Cells.Item(134, 11).Value = 20991.2901  # K134
$ws.Cells.Item(134, 12).Value = 14128.0005  # L134
$ws.Cells.Item(134, 13).Value = -18456.2901  # M134
$ws.Cells.Item(134, 14).Value = -19198.0005  # N134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2780.8708  # H31
$ws.Cells.Item(31, 9).Value = 1213.8096  # I31
$ws.Cells.Item(31, 10).Value = 6071.7  # J31
$ws.Cells.Item(31, 11).Value = 1213.8096  # K31
$ws.Cells.Item(31, 12).Value = 6071.7  # L31
$ws.Cells.Item(31, 13).Value = -918.8096  # M31
$ws.Cells.Item(31, 14).Value = -6661.7  # N31
$ws.Cells.Item(34, 8).Value = 2780.8708  # H34
$ws.Cells.Item(34, 9).Value = 1213.8096  # I34
$ws.Cells.Item(34, 10).Value = 6071.7  # J34
$ws.Cells.Item(34, 11).Value = 1213.8096  # K34
$ws.Cells.Item(34, 12).Value = 6071.7  # L34
$ws.Cells.Item(34, 13).Value = -1011.8096  # M34
$ws.Cells.Item(34, 14).Value = -6475.7  # N34
$ws.Cells.Item(37, 8).Value = 19800  # H37
$ws.Cells.Item(37, 10).Value = 19800  # J37
$ws.Cells.Item(37, 12).Value = 19800  # L37
$ws.Cells.Item(37, 14).Value = -20014  # N37
$ws.Cells.Item(45, 8).Value = 6698  # H45
$ws.Cells.Item(45, 9).Value = 6698  # I45
$ws.Cells.Item(45, 11).Value = 6698  # K45
$ws.Cells.Item(45, 13).Value = -6105  # M45
$ws.Cells.Item(58, 8).Value = 1188.8276  # H58
$ws.Cells.Item(58, 9).Value = 969.875  # I58
$ws.Cells.Item(58, 10).Value = 1458.3077  # J58
$ws.Cells.Item(58, 11).Value = 969.875  # K58
$ws.Cells.Item(58, 12).Value = 1458.3077  # L58
$ws.Cells.Item(58, 13).Value = -766.875  # M58
$ws.Cells.Item(58, 14).Value = -1864.3077  # N58
$ws.Cells.Item(106, 8).Value = 46566.332  # H106
$ws.Cells.Item(106, 10).Value = 46566.332  # J106
$ws.Cells.Item(106, 12).Value = 46566.332  # L106
$ws.Cells.Item(106, 14).Value = -49090.332  # N106
$ws.Cells.Item(122, 8).Value = 1382.1  # H122
$ws.Cells.Item(122, 9).Value = 1439.9474  # I122
$ws.Cells.Item(122, 10).Value = 1282.1818  # J122
$ws.Cells.Item(122, 11).Value = 4319.8422  # K122
$ws.Cells.Item(122, 12).Value = 3846.5454  # L122
$ws.Cells.Item(122, 13).Value = -1869.8422  # M122
$ws.Cells.Item(122, 14).Value = -8746.545399999999  # N122
$ws.Cells.Item(134, 8).Value = 965.25  # H134
$ws.Cells.Item(134, 9).Value = 917.5714  # I134
$ws.Cells.Item(134, 10).Value = 1299  # J134
$ws.Cells.Item(134, 11).Value = 2752.7142  # K134
$ws.Cells.Item(134, 12).Value = 3897  # L134
$ws.Cells.Item(134, 13).Value = -217.7142000000003  # M134
$ws.Cells.Item(134, 14).Value = -8967  # N134
$ws.Cells.Item(136, 8).Value = 1188.8276  # H136
$ws.Cells.Item(136, 9).Value = 969.875  # I136
$ws.Cells.Item(136, 10).Value = 1458.3077  # J136
$ws.Cells.Item(136, 11).Value = 2909.625  # K136
$ws.Cells.Item(136, 12).Value = 4374.9231  # L136
$ws.Cells.Item(136, 13).Value = -359.625  # M136
$ws.Cells.Item(136, 14).Value = -9474.9231  # N136

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 2577.6667  # H3
$ws.Cells.Item(3, 9).Value = 1200  # I3
$ws.Cells.Item(3, 10).Value = 3955.3333  # J3
$ws.Cells.Item(3, 11).Value = 3600  # K3
$ws.Cells.Item(3, 12).Value = 11865.9999  # L3
$ws.Cells.Item(3, 13).Value = -3488  # M3
$ws.Cells.Item(3, 14).Value = -12089.9999  # N3
$ws.Cells.Item(38, 8).Value = 273.06668  # H38
$ws.Cells.Item(38, 9).Value = 40.9  # I38
$ws.Cells.Item(38, 11).Value = 122.7  # K38
$ws.Cells.Item(38, 13).Value = 224.3  # M38
$ws.Cells.Item(68, 8).Value = 0  # H68
$ws.Cells.Item(68, 10).Value = 0  # J68
$ws.Cells.Item(68, 12).Value = 0  # L68
$ws.Cells.Item(68, 14).Value = ""  # N68: delete (was present, now cleared)
$ws.Cells.Item(71, 8).Value = 0  # H71
$ws.Cells.Item(71, 10).Value = 0  # J71
$ws.Cells.Item(71, 12).Value = 0  # L71
$ws.Cells.Item(71, 14).Value = ""  # N71: delete (was present, now cleared)
$ws.Cells.Item(87, 8).Value = 9654.333000000001  # H87
$ws.Cells.Item(87, 9).Value = 1981.5  # I87
$ws.Cells.Item(87, 11).Value = 5944.5  # K87
$ws.Cells.Item(87, 13).Value = -4696.5  # M87
$ws.Cells.Item(90, 8).Value = 9654.333000000001  # H90
$ws.Cells.Item(90, 9).Value = 1981.5  # I90
$ws.Cells.Item(90, 11).Value = 17833.5  # K90
$ws.Cells.Item(90, 13).Value = -11593.5  # M90
$ws.Cells.Item(129, 8).Value = 43566.707  # H129
$ws.Cells.Item(129, 10).Value = 61482.332  # J129
$ws.Cells.Item(129, 12).Value = 184446.996  # L129
$ws.Cells.Item(129, 14).Value = -194446.996  # N129
$ws.Cells.Item(131, 8).Value = 7153596.5  # H131
$ws.Cells.Item(131, 10).Value = 11691.141  # J131
$ws.Cells.Item(131, 12).Value = 35073.423  # L131
$ws.Cells.Item(131, 14).Value = -45153.423  # N131
$ws.Cells.Item(134, 8).Value = 1758.091  # H134
$ws.Cells.Item(134, 9).Value = 1404.4736  # I134
$ws.Cells.Item(134, 11).Value = 4213.4208  # K134
$ws.Cells.Item(134, 13).Value = 856.5792000000001  # M134
$ws.Cells.Item(137, 8).Value = 4973.8184  # H137
$ws.Cells.Item(137, 10).Value = 6924.5713  # J137
$ws.Cells.Item(137, 12).Value = 20773.7139  # L137
$ws.Cells.Item(137, 14).Value = -30973.7139  # N137
$ws.Cells.Item(139, 8).Value = 10099.667  # H139
$ws.Cells.Item(139, 9).Value = 11720  # I139
$ws.Cells.Item(139, 11).Value = 35160  # K139
$ws.Cells.Item(139, 13).Value = -30020  # M139

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 2205.5386  # H122
$ws.Cells.Item(122, 9).Value = 2170.5  # I122
$ws.Cells.Item(122, 10).Value = 2261.6  # J122
$ws.Cells.Item(122, 11).Value = 6511.5  # K122
$ws.Cells.Item(122, 12).Value = 6784.799999999999  # L122
$ws.Cells.Item(122, 13).Value = -4061.5  # M122
$ws.Cells.Item(122, 14).Value = -11684.8  # N122
$ws.Cells.Item(132, 8).Value = 3025.7  # H132
$ws.Cells.Item(132, 9).Value = 2534.9375  # I132
$ws.Cells.Item(132, 10).Value = 4988.75  # J132
$ws.Cells.Item(132, 11).Value = 7604.8125  # K132
$ws.Cells.Item(132, 12).Value = 14966.25  # L132
$ws.Cells.Item(132, 13).Value = -5074.8125  # M132
$ws.Cells.Item(132, 14).Value = -20026.25  # N132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 1937.88  # H132
$ws.Cells.Item(132, 9).Value = 1483.091  # I132
$ws.Cells.Item(132, 11).Value = 4449.272999999999  # K132
$ws.Cells.Item(132, 13).Value = -1919.272999999999  # M132

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 3702  # H62
$ws.Cells.Item(62, 9).Value = 3702  # I62
$ws.Cells.Item(62, 11).Value = 3702  # K62
$ws.Cells.Item(62, 13).Value = -3078  # M62
$ws.Cells.Item(65, 8).Value = 3702  # H65
$ws.Cells.Item(65, 9).Value = 3702  # I65
$ws.Cells.Item(65, 11).Value = 18510  # K65
$ws.Cells.Item(65, 13).Value = -15390  # M65
$ws.Cells.Item(80, 8).Value = 73999.5  # H80
$ws.Cells.Item(80, 10).Value = 73999.5  # J80
$ws.Cells.Item(80, 12).Value = 73999.5  # L80
$ws.Cells.Item(80, 14).Value = -75995.5  # N80
$ws.Cells.Item(83, 8).Value = 73999.5  # H83
$ws.Cells.Item(83, 10).Value = 73999.5  # J83
$ws.Cells.Item(83, 12).Value = 221998.5  # L83
$ws.Cells.Item(83, 14).Value = -231982.5  # N83
$ws.Cells.Item(111, 8).Value = 0  # H111
$ws.Cells.Item(111, 10).Value = 0  # J111
$ws.Cells.Item(111, 12).Value = 0  # L111
$ws.Cells.Item(111, 14).Value = ""  # N111: delete (was present, now cleared)
$ws.Cells.Item(122, 8).Value = 1473  # H122
$ws.Cells.Item(122, 9).Value = 1253.0588  # I122
$ws.Cells.Item(122, 11).Value = 3759.1764  # K122
$ws.Cells.Item(122, 13).Value = -1309.1764  # M122
$ws.Cells.Item(132, 8).Value = 4061.7307  # H132
$ws.Cells.Item(132, 9).Value = 2740.75  # I132
$ws.Cells.Item(132, 10).Value = 5194  # J132
$ws.Cells.Item(132, 11).Value = 8222.25  # K132
$ws.Cells.Item(132, 12).Value = 15582  # L132
$ws.Cells.Item(132, 13).Value = -5692.25  # M132
$ws.Cells.Item(132, 14).Value = -20642  # N132
